# Update cryptos list data (prices in column D, volume(1h) % in column E)
# Note: values that look numeric are prefixed with a leading apostrophe so
# Excel stores them as text (matching the original inline-string cells)
# rather than silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.735.95"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "1.640.77"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  +0.53%  "
$ws.Range("D5").Value = "'217.73"
$ws.Range("E5").Value = "  +1.40%  "
$ws.Range("D6").Value = "'0.503"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").Value = "'0.0625"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "'19.12"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "1.869.68"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "1.642.28"
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").Value = "'64.67"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("D17").Value = "26.731.09"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("E18").Value = "  -1.11%  "
$ws.Range("D19").Value = "'214.22"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").Value = "'2.38"
$ws.Range("E22").Value = "  +7.97%  "
$ws.Range("D23").Value = "'6.22"
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("D24").Value = "'9.27"
$ws.Range("E24").Value = "  -1.91%  "
$ws.Range("D25").Value = "'145.65"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").Value = "'0.118"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("D28").Value = "'7.16"
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("D29").Value = "'15.67"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("E31").Value = "  +1.40%  "
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("D34").Value = "1.287.24"
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("D35").Value = "'1.53"
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("D36").Value = "'2.43"
$ws.Range("E36").Value = "  +1.10%  "
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").Value = "'0.536"
$ws.Range("D39").Value = "'0.817"
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("E43").Value = "  -2.27%  "
$ws.Range("D44").Value = "1.779.12"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'60.98"
$ws.Range("E45").Value = "  +3.09%  "
$ws.Range("D46").Value = "'91.61"
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("D47").Value = "'1.59"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D49").Value = "'7.61"
$ws.Range("E49").Value = "  -1.28%  "
$ws.Range("D50").Value = "'0.0966"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("D51").Value = "'0.407"
$ws.Range("E51").Value = "  +0.35%  "
